$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cell H1 (bold, centered, bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new data columns I (I0) and J (IF) for rows 2-27.
$iValues = @(3,6,9,9,6,6,6,5,9,7,7,8,5,3,5,8,7,7,6,8,7,9,8,8,3,6)
$jValues = @(3,6,9,9,6,6,6,6,9,8,7,8,5,5,5,8,7,7,6,8,7,9,8,8,3,6)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
